# The lookup helper used to stash the data row's ordinal position in column C
# and read it back positionally; it now locates rows by an absolute cell
# reference instead, so that helper column value is no longer written.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data Sheet")

# Drop the now-unused ordinal-position value that used to live in C4.
$ws.Range("C4").ClearContents()

# Column A (the absolute reference column) is now the important one to read,
# so size it to fit its contents.
$ws.Columns.Item(1).EntireColumn.AutoFit()

# Leave the selection where the user's last action landed.
$ws.Range("A5").Select()
